$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data at the bottom: "Plaza Vea Universitaria" with factor 4
$ws.Range("A5").Value = "Plaza Vea Universitaria"
$ws.Range("B5").Value = 4
